{"js": "// Update the regression-results table:\n//  - \"BIC\" row values change from 9925.1/11052.8 to 5715.8/6232.9\n//  - a new \"F\" row (24.785 / 19.220) is inserted right after the BIC row\n//    (i.e. right before the trailing \"RMSE\" row).\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Find the BIC row by reading its first-cell text.\nconst rowItems = rows.items;\nfor (const row of rowItems) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\nlet bicRow = null;\nfor (const row of rowItems) {\n  const firstCellText = row.cells.items[0].value;\n  if (firstCellText && firstCellText.trim() === \"BIC\") {\n    bicRow = row;\n    break;\n  }\n}\n\nif (!bicRow) {\n  throw new Error(\"Could not locate the BIC row\");\n}\n\n// Update the BIC row's numeric values.\nbicRow.cells.items[1].value = \"5715.8\";\nbicRow.cells.items[2].value = \"6232.9\";\n\n// Insert a new row directly after the BIC row with the F-statistic values.\nbicRow.insertRows(\"After\", 1, [[\"F\", \"24.785\", \"19.220\"]]);\n\nawait context.sync();\n", "ps1": "# Update the regression-results table:\n#  - \"BIC\" row values change from 9925.1/11052.8 to 5715.8/6232.9\n#  - a new \"F\" row (24.785 / 19.220) is inserted right after the BIC row\n#    (i.e. right before the trailing \"RMSE\" row).\n\n$d = $word.ActiveDocument\n$t = $d.Tables(1)\n\n# Locate the BIC row by scanning the first column. Cell.Range.Text carries a\n# trailing cell-mark (CR + BEL) that needs stripping before comparing.\n$bicRow = $null\n$rmseRow = $null\nfor ($i = 1; $i -le $t.Rows.Count; $i++) {\n    $label = $t.Rows($i).Cells(1).Range.Text.TrimEnd([char]13, [char]7).Trim()\n    if ($label -eq \"BIC\") {\n        $bicRow = $t.Rows($i)\n    }\n    if ($label -eq \"RMSE\") {\n        $rmseRow = $t.Rows($i)\n    }\n}\n\nif ($bicRow -eq $null) {\n    throw \"Could not locate the BIC row\"\n}\n\n# Update the BIC row's numeric values.\n$bicRow.Cells(2).Range.Text = \"5715.8\"\n$bicRow.Cells(3).Range.Text = \"6232.9\"\n\n# Insert a new row directly before the trailing RMSE row (i.e. right after BIC)\n# and populate it with the F-statistic values.\nif ($rmseRow -ne $null) {\n    $newRow = $t.Rows.Add($rmseRow)\n} else {\n    $newRow = $t.Rows.Add()\n}\n$newRow.Cells(1).Range.Text = \"F\"\n$newRow.Cells(2).Range.Text = \"24.785\"\n$newRow.Cells(3).Range.Text = \"19.220\"\n"}
